$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'92.390.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.18%  "

$ws.Range("D3").Value = "'3.317.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.11%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'227.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -10.57%  "

$ws.Range("D6").Value = "'620.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.97%  "

$ws.Range("D7").Value = "'1.33"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -10.93%  "

$ws.Range("D8").Value = "'0.376"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -12.32%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("D10").Value = "'0.910"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -13.79%  "

$ws.Range("D11").Value = "'3.306.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.34%  "

$ws.Range("D12").Value = "'0.191"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.74%  "

$ws.Range("D13").Value = "'39.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -13.60%  "

$ws.Range("D14").Value = "'92.230.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.10%  "

$ws.Range("D15").Value = "'5.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.91%  "

$ws.Range("D16").Value = "'3.926.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.52%  "

$ws.Range("D17").Value = "'0.0000240"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -9.06%  "

$ws.Range("D18").Value = "'7.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -13.67%  "

$ws.Range("D19").Value = "'3.325.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.19%  "

$ws.Range("D20").Value = "'16.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -12.37%  "

$ws.Range("D21").Value = "'10.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.56%  "

$ws.Range("D22").Value = "'482.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.19%  "

$ws.Range("D23").Value = "'0.438"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -16.85%  "

$ws.Range("D24").Value = "'3.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -11.11%  "

$ws.Range("D25").Value = "'0.0000181"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -11.27%  "

$ws.Range("D26").Value = "'6.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -11.58%  "

$ws.Range("D27").Value = "'88.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.57%  "

$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'11.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -11.93%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'11.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.03%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.24%  "

$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.129"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -11.66%  "

$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D34").Value = "'0.167"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -12.38%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'27.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -10.38%  "

$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "'0.515"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -14.47%  "

$ws.Range("B37").Value = "USDe"
$ws.Range("C37").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'7.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.35%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'510.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.40%  "

$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "'1.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.98%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.144"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.49%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'0.857"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.62%  "

$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'23.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.94%  "

$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").Value = "'1.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.66%  "

$ws.Range("D45").Value = "'3.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.32%  "

$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'5.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.31%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'2.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.83%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0385"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.80%  "

$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'51.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.23%  "

$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "'3.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.09%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'7.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.27%  "
